$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture "before" values for the rows involved in the cyclic rotation.
$rows = @(2, 3, 5, 6, 7, 8, 9)
$cols = @("D", "K", "L", "M", "N", "O", "P", "R", "S")

$before = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowData
}

# Mapping: new value of row X = old value of row Map[X]
$map = @{
    2 = 3
    3 = 9
    9 = 7
    7 = 5
    5 = 6
    6 = 8
    8 = 2
}

foreach ($r in $rows) {
    $src = $map[$r]
    $srcData = $before[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $srcData[$c]
    }
}
